$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting Name/Age/Grade/Score (and
# the student rows) one column to the right (A:D -> B:E).
$ws.Range("A1").EntireColumn.Insert()

# Give the new column A's data rows the same look as the header row (bold,
# bordered, centered) by copying the header cell's formatting only.
$ws.Range("B1").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)

# Fill in column A with a 0-based row index for each student.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
